$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2931.5356
$ws.Range("I107").Value = 2279.182
$ws.Range("J107").Value = 5323.5
$ws.Range("K107").Value = 2279.182
$ws.Range("L107").Value = 5323.5
$ws.Range("M107").Value = -359.1819999999998
$ws.Range("N107").Value = -9163.5
$ws.Range("H112").Value = 1854.7142
$ws.Range("J112").Value = 2041.7646
$ws.Range("L112").Value = 6125.293799999999
$ws.Range("N112").Value = -8341.293799999999
$ws.Range("H129").Value = 853.86365
$ws.Range("I129").Value = 694
$ws.Range("J129").Value = 889.3889
$ws.Range("K129").Value = 2082
$ws.Range("L129").Value = 2668.1667
$ws.Range("M129").Value = 2918
$ws.Range("N129").Value = -12668.1667
$ws.Range("H132").Value = 7099315.5
$ws.Range("I132").Value = 8134686
$ws.Range("K132").Value = 24404058
$ws.Range("M132").Value = -24401528
$ws.Range("H137").Value = 2651.0889
$ws.Range("I137").Value = 2264.3
$ws.Range("K137").Value = 6792.900000000001
$ws.Range("M137").Value = -4242.900000000001
$ws.Range("H138").Value = 2553.5
$ws.Range("I138").Value = 2315.6667
$ws.Range("J138").Value = 2588.305
$ws.Range("K138").Value = 6947.000100000001
$ws.Range("L138").Value = 7764.914999999999
$ws.Range("M138").Value = -1807.000100000001
$ws.Range("N138").Value = -18044.915

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9805.904
$ws.Range("I32").Value = 7099.364
$ws.Range("J32").Value = 22064.941
$ws.Range("K32").Value = 7099.364
$ws.Range("L32").Value = 22064.941
$ws.Range("M32").Value = -6812.364
$ws.Range("N32").Value = -22638.941
$ws.Range("H61").Value = 142859090
$ws.Range("I61").Value = 200001520
$ws.Range("J61").Value = 2999.5
$ws.Range("K61").Value = 200001520
$ws.Range("L61").Value = 2999.5
$ws.Range("M61").Value = -200001308
$ws.Range("N61").Value = -3423.5
$ws.Range("H80").Value = 37750
$ws.Range("J80").Value = 37750
$ws.Range("L80").Value = 37750
$ws.Range("N80").Value = -39746
$ws.Range("H83").Value = 37750
$ws.Range("J83").Value = 37750
$ws.Range("L83").Value = 113250
$ws.Range("N83").Value = -123234
$ws.Range("H110").Value = 253.8
$ws.Range("I110").Value = 258
$ws.Range("J110").Value = 237
$ws.Range("K110").Value = 258
$ws.Range("L110").Value = 237
$ws.Range("M110").Value = 1787
$ws.Range("N110").Value = -4327
$ws.Range("H112").Value = 48333.332
$ws.Range("J112").Value = 48333.332
$ws.Range("L112").Value = 48333.332
$ws.Range("N112").Value = -51287.332
$ws.Range("H132").Value = 3937.8462
$ws.Range("I132").Value = 3887.3333
$ws.Range("J132").Value = 4051.5
$ws.Range("K132").Value = 11661.9999
$ws.Range("L132").Value = 12154.5
$ws.Range("M132").Value = -9131.999899999999
$ws.Range("N132").Value = -17214.5
$ws.Range("H134").Value = 36260
$ws.Range("J134").Value = 36260
$ws.Range("L134").Value = 36260
$ws.Range("N134").Value = -46400
$ws.Range("H136").Value = 142859090
$ws.Range("I136").Value = 200001520
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 600004560
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -600002010
$ws.Range("N136").Value = -14098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 10009
$ws.Range("J14").Value = 10009
$ws.Range("L14").Value = 10009
$ws.Range("N14").Value = -10353
$ws.Range("H86").Value = 4465.643
$ws.Range("I86").Value = 4914.125
$ws.Range("J86").Value = 3867.6667
$ws.Range("K86").Value = 4914.125
$ws.Range("L86").Value = 3867.6667
$ws.Range("M86").Value = -3791.125
$ws.Range("N86").Value = -6113.6667
$ws.Range("H89").Value = 4465.643
$ws.Range("I89").Value = 4914.125
$ws.Range("J89").Value = 3867.6667
$ws.Range("K89").Value = 24570.625
$ws.Range("L89").Value = 19338.3335
$ws.Range("M89").Value = -18954.625
$ws.Range("N89").Value = -30570.3335
$ws.Range("H110").Value = 30333
$ws.Range("J110").Value = 30333
$ws.Range("L110").Value = 30333
$ws.Range("N110").Value = -38513
$ws.Range("H134").Value = 3957.4546
$ws.Range("I134").Value = 929.89655
$ws.Range("K134").Value = 2789.68965
$ws.Range("M134").Value = -254.6896500000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 83334960
$ws.Range("I16").Value = 90910720
$ws.Range("J16").Value = 1513
$ws.Range("K16").Value = 90910720
$ws.Range("L16").Value = 1513
$ws.Range("M16").Value = -90910433
$ws.Range("N16").Value = -2087
$ws.Range("H31").Value = 1642.4791
$ws.Range("I31").Value = 1494.8422
$ws.Range("K31").Value = 1494.8422
$ws.Range("M31").Value = -1199.8422
$ws.Range("H34").Value = 1642.4791
$ws.Range("I34").Value = 1494.8422
$ws.Range("K34").Value = 1494.8422
$ws.Range("M34").Value = -1292.8422
$ws.Range("H58").Value = 14512.375
$ws.Range("I58").Value = 2120
$ws.Range("J58").Value = 35166.332
$ws.Range("K58").Value = 2120
$ws.Range("L58").Value = 35166.332
$ws.Range("M58").Value = -1917
$ws.Range("N58").Value = -35572.332
$ws.Range("H107").Value = 652.6923
$ws.Range("I107").Value = 323.09525
$ws.Range("K107").Value = 323.09525
$ws.Range("M107").Value = 1596.90475
$ws.Range("H113").Value = 83334960
$ws.Range("I113").Value = 90910720
$ws.Range("J113").Value = 1513
$ws.Range("K113").Value = 90910720
$ws.Range("L113").Value = 1513
$ws.Range("M113").Value = -90908550
$ws.Range("N113").Value = -5853
$ws.Range("H132").Value = 2056.913
$ws.Range("J132").Value = 2762.6
$ws.Range("L132").Value = 8287.799999999999
$ws.Range("N132").Value = -13347.8
$ws.Range("H134").Value = 15626380
$ws.Range("I134").Value = 1345.375
$ws.Range("K134").Value = 4036.125
$ws.Range("M134").Value = -1501.125
$ws.Range("H136").Value = 14512.375
$ws.Range("I136").Value = 2120
$ws.Range("J136").Value = 35166.332
$ws.Range("K136").Value = 6360
$ws.Range("L136").Value = 105498.996
$ws.Range("M136").Value = -3810
$ws.Range("N136").Value = -110598.996
$ws.Range("H141").Value = 323620.8
$ws.Range("J141").Value = 323620.8
$ws.Range("L141").Value = 323620.8
$ws.Range("N141").Value = -333980.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 696.89795
$ws.Range("I113").Value = 620.1579
$ws.Range("J113").Value = 745.5
$ws.Range("K113").Value = 1860.4737
$ws.Range("L113").Value = 2236.5
$ws.Range("M113").Value = 309.5263
$ws.Range("N113").Value = -6576.5
$ws.Range("H131").Value = 20865568
$ws.Range("I131").Value = 250000450
$ws.Range("J131").Value = 35124.297
$ws.Range("K131").Value = 750001350
$ws.Range("L131").Value = 105372.891
$ws.Range("M131").Value = -749996310
$ws.Range("N131").Value = -115452.891

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 40912716
$ws.Range("I70").Value = 35718044
$ws.Range("J70").Value = 50003396
$ws.Range("K70").Value = 35718044
$ws.Range("L70").Value = 50003396
$ws.Range("M70").Value = -35717774
$ws.Range("N70").Value = -50003936
$ws.Range("H73").Value = 40912716
$ws.Range("I73").Value = 35718044
$ws.Range("J73").Value = 50003396
$ws.Range("K73").Value = 35718044
$ws.Range("L73").Value = 50003396
$ws.Range("M73").Value = -35717108
$ws.Range("N73").Value = -50005268
$ws.Range("H113").Value = 1378.8889
$ws.Range("I113").Value = 1730
$ws.Range("K113").Value = 1730
$ws.Range("M113").Value = 440
$ws.Range("H132").Value = 4927.6216
$ws.Range("I132").Value = 5840.04
$ws.Range("J132").Value = 3026.75
$ws.Range("K132").Value = 17520.12
$ws.Range("L132").Value = 9080.25
$ws.Range("M132").Value = -14990.12
$ws.Range("N132").Value = -14140.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 742.5185
$ws.Range("I16").Value = 767.61536
$ws.Range("J16").Value = 90
$ws.Range("K16").Value = 767.61536
$ws.Range("L16").Value = 90
$ws.Range("M16").Value = -597.61536
$ws.Range("N16").Value = -430
$ws.Range("H82").Value = 2316.9
$ws.Range("I82").Value = 2271.125
$ws.Range("K82").Value = 2271.125
$ws.Range("M82").Value = -1910.125
$ws.Range("H85").Value = 2316.9
$ws.Range("I85").Value = 2271.125
$ws.Range("K85").Value = 2271.125
$ws.Range("M85").Value = -1023.125
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180
$ws.Range("H132").Value = 2212.761
$ws.Range("I132").Value = 1868.7587
$ws.Range("J132").Value = 2799.5881
$ws.Range("K132").Value = 5606.2761
$ws.Range("L132").Value = 8398.764299999999
$ws.Range("M132").Value = -3076.2761
$ws.Range("N132").Value = -13458.7643

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 556.55554
$ws.Range("I107").Value = 501.125
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1503.375
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 416.625
$ws.Range("N107").Value = -6840
$ws.Range("H132").Value = 2144.7446
$ws.Range("I132").Value = 2127.639
$ws.Range("K132").Value = 6382.917
$ws.Range("M132").Value = -3852.917
$ws.Range("H141").Value = 53950.715
$ws.Range("J141").Value = 53950.715
$ws.Range("L141").Value = 53950.715
$ws.Range("N141").Value = -64310.715
